$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "87.840.56"
$ws.Range("E2").Value = "  +10.76%  "
$ws.Range("D3").Value = "3.351.76"
$ws.Range("E3").Value = "  +7.22%  "
$ws.Range("E4").Value = "  -0.04%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "220.16"
$c.ClearFormats()
$ws.Range("E5").Value = "  +7.63%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "640.51"
$c.ClearFormats()
$ws.Range("E6").Value = "  +3.48%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.323"
$c.ClearFormats()
$ws.Range("E7").Value = "  +23.69%  "
$ws.Range("E8").Value = "  -0.09%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.618"
$c.ClearFormats()
$ws.Range("E9").Value = "  +6.72%  "
$ws.Range("D10").Value = "3.357.54"
$ws.Range("E10").Value = "  +7.43%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.612"
$c.ClearFormats()
$ws.Range("E11").Value = "  +6.00%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.0000277"
$c.ClearFormats()
$ws.Range("E12").Value = "  +13.49%  "
$ws.Range("E13").Value = "  +2.62%  "
$ws.Range("D14").Value = "3.962.44"
$ws.Range("E14").Value = "  +7.25%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "34.55"
$c.ClearFormats()
$ws.Range("E15").Value = "  +11.13%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "5.42"
$c.ClearFormats()
$ws.Range("E16").Value = "  +4.94%  "
$ws.Range("D17").Value = "87.435.90"
$ws.Range("E17").Value = "  +10.47%  "
$ws.Range("D18").Value = "3.344.94"
$ws.Range("E18").Value = "  +7.72%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "14.72"
$c.ClearFormats()
$ws.Range("E19").Value = "  +5.04%  "
$ws.Range("B20").Value = "SuiNetwork"
$ws.Range("C20").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "3.23"
$c.ClearFormats()
$ws.Range("E20").Value = "  +10.94%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "451.19"
$c.ClearFormats()
$ws.Range("E21").Value = "  +5.17%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "9.21"
$c.ClearFormats()
$ws.Range("E22").Value = "  +2.74%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "5.34"
$c.ClearFormats()
$ws.Range("E23").Value = "  +3.96%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "7.40"
$c.ClearFormats()
$ws.Range("E24").Value = "  +9.15%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "5.38"
$c.ClearFormats()
$ws.Range("E25").Value = "  +16.75%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "12.37"
$c.ClearFormats()
$ws.Range("E26").Value = "  +16.02%  "
$ws.Range("D27").Value = "3.499.71"
$ws.Range("E27").Value = "  +6.50%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "79.05"
$c.ClearFormats()
$ws.Range("E28").Value = "  +5.50%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0.0000131"
$c.ClearFormats()
$ws.Range("E29").Value = "  +10.44%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.ClearFormats()
$ws.Range("E30").Value = "  +0.03%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.189"
$c.ClearFormats()
$ws.Range("E31").Value = "  +55.59%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "604.87"
$c.ClearFormats()
$ws.Range("E32").Value = "  +10.47%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "9.34"
$c.ClearFormats()
$ws.Range("E33").Value = "  +5.92%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.ClearFormats()
$ws.Range("E34").Value = "  -0.07%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.56"
$c.ClearFormats()
$ws.Range("E35").Value = "  +8.03%  "
$ws.Range("E36").Value = "  +5.08%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.153"
$c.ClearFormats()
$ws.Range("E37").Value = "  +2.39%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "23.57"
$c.ClearFormats()
$ws.Range("E38").Value = "  +4.42%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "6.64"
$c.ClearFormats()
$ws.Range("E39").Value = "  +21.37%  "
$ws.Range("E40").Value = "  +6.16%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.ClearFormats()
$ws.Range("E41").Value = "  +0.01%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "21.39"
$c.ClearFormats()
$ws.Range("E42").Value = "  +3.25%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "3.13"
$c.ClearFormats()
$ws.Range("E43").Value = "  +18.61%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "2.07"
$c.ClearFormats()
$ws.Range("E44").Value = "  +16.67%  "
$ws.Range("E45").Value = "  +0.11%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "157.41"
$c.ClearFormats()
$ws.Range("E46").Value = "  -2.98%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "190.90"
$c.ClearFormats()
$ws.Range("E47").Value = "  +2.34%  "
$ws.Range("E48").Value = "  +8.49%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "45.95"
$c.ClearFormats()
$ws.Range("E49").Value = "  +8.15%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.792"
$c.ClearFormats()
$ws.Range("E50").Value = "  +3.64%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "26.64"
$c.ClearFormats()
$ws.Range("E51").Value = "  +10.33%  "
